$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("M2")
$r.Borders.ColorIndex = -4105
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1
